$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "245×2=490"
$t.Cell(1,2).Range.Text = "603×2=1206"
$t.Cell(1,3).Range.Text = "222×2=444"
$t.Cell(1,4).Range.Text = "874×9=7866"
$t.Cell(1,5).Range.Text = "760×6=4560"
$t.Cell(5,1).Range.Text = "309×4=1236"
$t.Cell(5,2).Range.Text = "179×6=1074"
$t.Cell(5,3).Range.Text = "353×6=2118"
$t.Cell(5,4).Range.Text = "148×3=444"
$t.Cell(5,5).Range.Text = "805×7=5635"
$t.Cell(10,1).Range.Text = "878×8=7024"
$t.Cell(10,2).Range.Text = "367×4=1468"
$t.Cell(10,3).Range.Text = "314×4=1256"
$t.Cell(10,4).Range.Text = "523×2=1046"
$t.Cell(10,5).Range.Text = "350×3=1050"
$t.Cell(15,1).Range.Text = "678×3=2034"
$t.Cell(15,2).Range.Text = "132×3=396"
$t.Cell(15,3).Range.Text = "608×8=4864"
$t.Cell(15,4).Range.Text = "575×9=5175"
$t.Cell(15,5).Range.Text = "634×2=1268"
$t.Cell(20,1).Range.Text = "528×3=1584"
$t.Cell(20,2).Range.Text = "848×9=7632"
$t.Cell(20,3).Range.Text = "382×5=1910"
$t.Cell(20,4).Range.Text = "257×6=1542"
$t.Cell(20,5).Range.Text = "985×6=5910"
